# EventPatterns.xlsx - add two new event pattern rows (unnamed-group regex
# support) to the TopPatterns table, and re-point the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TopPatterns")

# --- Insert "OldLogFilesMoved" above the existing "XMLTransactionFailure"
#     row (old row 6), pushing it and everything below down by one. ---
$ws.Rows.Item(6).Insert()
$ws.Range("E6").Value = "OldLogFilesMoved"
$ws.Range("F6").Value = "Old log files moved (into|to)\s*(?P<_folder>.*)"

# --- Insert "NoBackendsLoaded" above the existing "AlreadyLoggedIn" row
#     (old row 7, now row 8 after the first insert), pushing it and
#     everything below down by one. ---
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 6
$ws.Range("E8").Value = "NoBackendsLoaded"
$ws.Range("F8").Value = "(?P<user>.*): No backends loaded.(?P<_errorDetail>.*)"

# Table1_3 (the ListObject bound to A1:G8) auto-extends to A1:G10 on the
# row inserts above, matching the widened sheet dimension.

# Drop the stale row-outline level left on the sheet (no rows are actually
# grouped) while keeping the 6-level column outline intact.
$ws.Outline.ShowLevels(0, 6)

# Move the active selection to reflect where the edit left off.
$ws.Range("F7").Select() | Out-Null
